$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string must be forced to Text
# format first, otherwise Excel auto-converts the assigned string to a number
# (the source data models these Price cells as text, e.g. "569.56").
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'

$ws.Range('D2').Value = '60.841.58'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '3.366.43'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '569.56'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').Value = '138.82'
$ws.Range('E6').Value = '  -2.01%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('D9').Value = '7.65'
$ws.Range('E9').Value = '  +1.78%  '
$ws.Range('E10').Value = '  -2.42%  '
$ws.Range('E11').Value = '  -4.54%  '
$ws.Range('D12').Value = '3.941.82'
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').Value = '0.126'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('D14').Value = '27.71'
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('D15').Value = '3.373.92'
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('E16').Value = '  -2.04%  '
$ws.Range('D17').Value = '60.931.41'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('E18').Value = '  -3.14%  '
$ws.Range('D19').Value = '13.54'
$ws.Range('E19').Value = '  -3.57%  '
$ws.Range('D20').Value = '8.88'
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('D21').Value = '381.68'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').Value = '75.32'
$ws.Range('E22').Value = '  +2.16%  '
$ws.Range('D23').Value = '0.548'
$ws.Range('E23').Value = '  -2.47%  '
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('E25').Value = '  -5.50%  '
$ws.Range('E26').Value = '  +6.42%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  -3.80%  '
$ws.Range('D29').Value = '7.81'
$ws.Range('E29').Value = '  -2.26%  '
$ws.Range('E30').Value = '  -2.12%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  -7.06%  '
$ws.Range('D33').Value = '22.91'
$ws.Range('E33').Value = '  -3.34%  '
$ws.Range('E34').Value = '  -1.99%  '
$ws.Range('D35').Value = '166.25'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').Value = '4.91'
$ws.Range('E36').Value = '  -2.04%  '
$ws.Range('D37').Value = '3.402.98'
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('E38').Value = '  -3.84%  '
$ws.Range('E39').Value = '  -2.33%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Value = '25.25'
$ws.Range('E40').Value = '  -9.39%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').Value = '0.773'
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('E42').Value = '  -2.74%  '
$ws.Range('E43').Value = '  -3.79%  '
$ws.Range('E44').Value = '  -1.92%  '
$ws.Range('D45').Value = '2.447.56'
$ws.Range('E45').Value = '  -3.66%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('E47').Value = '  -4.00%  '
$ws.Range('D48').Value = '22.16'
$ws.Range('E48').Value = '  -5.75%  '
$ws.Range('E49').Value = '  -4.76%  '
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('E51').Value = '  -3.60%  '
